$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the dataset (RM 232, SC 92).
# Deleting row 26 first shifts everything below it up by one, so the row
# that used to be "SC 92" (originally row 28) is now at row 27.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Re-impute / clear individual cells so the remaining rows match the
# updated (error-recalculated) values.
$ws.Range("C2").Value = 14.9

$ws.Range("C3").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("E5").Value = ""

$ws.Range("E8").Value = -6.6
$ws.Range("E10").Value = -6.1
$ws.Range("C11").Value = 11.4

$ws.Range("E12").Value = ""
$ws.Range("C13").Value = ""

$ws.Range("E15").Value = -8.4

$ws.Range("E18").Value = ""
$ws.Range("E19").Value = ""

$ws.Range("C21").Value = 12.7

$ws.Range("C25").Value = ""
$ws.Range("E25").Value = -7.1

$ws.Range("E27").Value = -10

$ws.Range("B29").Value = ""
$ws.Range("E29").Value = ""

$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("E33").Value = ""

Write-Output "done"
